$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Refresh the time_taken column (F2:F401) on the "data" sheet ---
# (regenerated values from re-running the panelapp query script)
$newTimeTaken = @(
  "2021-10-05 14:35:30.563852",
  "2021-10-05 14:35:30.563860",
  "2021-10-05 14:35:30.563863",
  "2021-10-05 14:35:30.563865",
  "2021-10-05 14:35:30.563868",
  "2021-10-05 14:35:30.563871",
  "2021-10-05 14:35:30.563873",
  "2021-10-05 14:35:30.563876",
  "2021-10-05 14:35:30.563878",
  "2021-10-05 14:35:30.563881",
  "2021-10-05 14:35:30.563883",
  "2021-10-05 14:35:30.563886",
  "2021-10-05 14:35:30.563888",
  "2021-10-05 14:35:30.563891",
  "2021-10-05 14:35:30.563893",
  "2021-10-05 14:35:30.563896",
  "2021-10-05 14:35:30.563898",
  "2021-10-05 14:35:30.563901",
  "2021-10-05 14:35:30.563904",
  "2021-10-05 14:35:30.563906",
  "2021-10-05 14:35:30.563909",
  "2021-10-05 14:35:30.563911",
  "2021-10-05 14:35:30.563914",
  "2021-10-05 14:35:30.563916",
  "2021-10-05 14:35:30.563919",
  "2021-10-05 14:35:30.563922",
  "2021-10-05 14:35:30.563924",
  "2021-10-05 14:35:30.563927",
  "2021-10-05 14:35:30.563929",
  "2021-10-05 14:35:30.563931",
  "2021-10-05 14:35:30.563934",
  "2021-10-05 14:35:30.563936",
  "2021-10-05 14:35:30.563939",
  "2021-10-05 14:35:30.563942",
  "2021-10-05 14:35:30.563944",
  "2021-10-05 14:35:30.563947",
  "2021-10-05 14:35:30.563949",
  "2021-10-05 14:35:30.563952",
  "2021-10-05 14:35:30.563954",
  "2021-10-05 14:35:30.563957",
  "2021-10-05 14:35:30.563959",
  "2021-10-05 14:35:30.563962",
  "2021-10-05 14:35:30.563965",
  "2021-10-05 14:35:30.563967",
  "2021-10-05 14:35:30.563970",
  "2021-10-05 14:35:30.563972",
  "2021-10-05 14:35:30.563974",
  "2021-10-05 14:35:30.563977",
  "2021-10-05 14:35:30.563979",
  "2021-10-05 14:35:30.563982",
  "2021-10-05 14:35:30.563984",
  "2021-10-05 14:35:30.563987",
  "2021-10-05 14:35:30.563990",
  "2021-10-05 14:35:30.563992",
  "2021-10-05 14:35:30.563995",
  "2021-10-05 14:35:30.563997",
  "2021-10-05 14:35:30.564000",
  "2021-10-05 14:35:30.564002",
  "2021-10-05 14:35:30.564005",
  "2021-10-05 14:35:30.564007",
  "2021-10-05 14:35:30.564010",
  "2021-10-05 14:35:30.564012",
  "2021-10-05 14:35:30.564014",
  "2021-10-05 14:35:30.564017",
  "2021-10-05 14:35:30.564020",
  "2021-10-05 14:35:30.564023",
  "2021-10-05 14:35:30.564026",
  "2021-10-05 14:35:30.564028",
  "2021-10-05 14:35:30.564031",
  "2021-10-05 14:35:30.564033",
  "2021-10-05 14:35:30.564035",
  "2021-10-05 14:35:30.564038",
  "2021-10-05 14:35:30.564040",
  "2021-10-05 14:35:30.564043",
  "2021-10-05 14:35:30.564045",
  "2021-10-05 14:35:30.564048",
  "2021-10-05 14:35:30.564052",
  "2021-10-05 14:35:30.564055",
  "2021-10-05 14:35:30.564058",
  "2021-10-05 14:35:30.564061",
  "2021-10-05 14:35:30.564063",
  "2021-10-05 14:35:30.564066",
  "2021-10-05 14:35:30.564068",
  "2021-10-05 14:35:30.564071",
  "2021-10-05 14:35:30.564073",
  "2021-10-05 14:35:30.564076",
  "2021-10-05 14:35:30.564078",
  "2021-10-05 14:35:30.564081",
  "2021-10-05 14:35:30.564083",
  "2021-10-05 14:35:30.564086",
  "2021-10-05 14:35:30.564088",
  "2021-10-05 14:35:30.564091",
  "2021-10-05 14:35:30.564094",
  "2021-10-05 14:35:30.564097",
  "2021-10-05 14:35:30.564100",
  "2021-10-05 14:35:30.564102",
  "2021-10-05 14:35:30.564105",
  "2021-10-05 14:35:30.564107",
  "2021-10-05 14:35:30.564110",
  "2021-10-05 14:35:30.564112",
  "2021-10-05 14:35:30.564115",
  "2021-10-05 14:35:30.564117",
  "2021-10-05 14:35:30.564120",
  "2021-10-05 14:35:30.564122",
  "2021-10-05 14:35:30.564125",
  "2021-10-05 14:35:30.564127",
  "2021-10-05 14:35:30.564130",
  "2021-10-05 14:35:30.564132",
  "2021-10-05 14:35:30.564136",
  "2021-10-05 14:35:30.564139",
  "2021-10-05 14:35:30.564142",
  "2021-10-05 14:35:30.564144",
  "2021-10-05 14:35:30.564146",
  "2021-10-05 14:35:30.564149",
  "2021-10-05 14:35:30.564151",
  "2021-10-05 14:35:30.564154",
  "2021-10-05 14:35:30.564156",
  "2021-10-05 14:35:30.564159",
  "2021-10-05 14:35:30.564161",
  "2021-10-05 14:35:30.564164",
  "2021-10-05 14:35:30.564166",
  "2021-10-05 14:35:30.564169",
  "2021-10-05 14:35:30.564171",
  "2021-10-05 14:35:30.564174",
  "2021-10-05 14:35:30.564176",
  "2021-10-05 14:35:30.564179",
  "2021-10-05 14:35:30.564181",
  "2021-10-05 14:35:30.564184",
  "2021-10-05 14:35:30.564187",
  "2021-10-05 14:35:30.564190",
  "2021-10-05 14:35:30.564193",
  "2021-10-05 14:35:30.564195",
  "2021-10-05 14:35:30.564198",
  "2021-10-05 14:35:30.564200",
  "2021-10-05 14:35:30.564203",
  "2021-10-05 14:35:30.564205",
  "2021-10-05 14:35:30.564208",
  "2021-10-05 14:35:30.564210",
  "2021-10-05 14:35:30.564213",
  "2021-10-05 14:35:30.564215",
  "2021-10-05 14:35:30.564218",
  "2021-10-05 14:35:30.564220",
  "2021-10-05 14:35:30.564223",
  "2021-10-05 14:35:30.564225",
  "2021-10-05 14:35:30.564228",
  "2021-10-05 14:35:30.564230",
  "2021-10-05 14:35:30.564233",
  "2021-10-05 14:35:30.564235",
  "2021-10-05 14:35:30.564238",
  "2021-10-05 14:35:30.564241",
  "2021-10-05 14:35:30.564243",
  "2021-10-05 14:35:30.564245",
  "2021-10-05 14:35:30.564248",
  "2021-10-05 14:35:30.564250",
  "2021-10-05 14:35:30.564253",
  "2021-10-05 14:35:30.564255",
  "2021-10-05 14:35:30.564258",
  "2021-10-05 14:35:30.564260",
  "2021-10-05 14:35:30.564263",
  "2021-10-05 14:35:30.564265",
  "2021-10-05 14:35:30.564268",
  "2021-10-05 14:35:30.564270",
  "2021-10-05 14:35:30.564273",
  "2021-10-05 14:35:30.564275",
  "2021-10-05 14:35:30.564278",
  "2021-10-05 14:35:30.564280",
  "2021-10-05 14:35:30.564282",
  "2021-10-05 14:35:30.564285",
  "2021-10-05 14:35:30.564287",
  "2021-10-05 14:35:30.564290",
  "2021-10-05 14:35:30.564292",
  "2021-10-05 14:35:30.564295",
  "2021-10-05 14:35:30.564299",
  "2021-10-05 14:35:30.564302",
  "2021-10-05 14:35:30.564304",
  "2021-10-05 14:35:30.564307",
  "2021-10-05 14:35:30.564309",
  "2021-10-05 14:35:30.564312",
  "2021-10-05 14:35:30.564314",
  "2021-10-05 14:35:30.564317",
  "2021-10-05 14:35:30.564319",
  "2021-10-05 14:35:30.564322",
  "2021-10-05 14:35:30.564324",
  "2021-10-05 14:35:30.564327",
  "2021-10-05 14:35:30.564329",
  "2021-10-05 14:35:30.564332",
  "2021-10-05 14:35:30.564334",
  "2021-10-05 14:35:30.564337",
  "2021-10-05 14:35:30.564339",
  "2021-10-05 14:35:30.564342",
  "2021-10-05 14:35:30.564344",
  "2021-10-05 14:35:30.564347",
  "2021-10-05 14:35:30.564349",
  "2021-10-05 14:35:30.564352",
  "2021-10-05 14:35:30.564354",
  "2021-10-05 14:35:30.564357",
  "2021-10-05 14:35:30.564359",
  "2021-10-05 14:35:30.564362",
  "2021-10-05 14:35:30.564364",
  "2021-10-05 14:35:30.564367",
  "2021-10-05 14:35:30.564370",
  "2021-10-05 14:35:30.564372",
  "2021-10-05 14:35:30.564375",
  "2021-10-05 14:35:30.564377",
  "2021-10-05 14:35:30.564380",
  "2021-10-05 14:35:30.564382",
  "2021-10-05 14:35:30.564385",
  "2021-10-05 14:35:30.564387",
  "2021-10-05 14:35:30.564390",
  "2021-10-05 14:35:30.564392",
  "2021-10-05 14:35:30.564395",
  "2021-10-05 14:35:30.564397",
  "2021-10-05 14:35:30.564400",
  "2021-10-05 14:35:30.564402",
  "2021-10-05 14:35:30.564405",
  "2021-10-05 14:35:30.564407",
  "2021-10-05 14:35:30.564410",
  "2021-10-05 14:35:30.564412",
  "2021-10-05 14:35:30.564415",
  "2021-10-05 14:35:30.564417",
  "2021-10-05 14:35:30.564419",
  "2021-10-05 14:35:30.564422",
  "2021-10-05 14:35:30.564424",
  "2021-10-05 14:35:30.564427",
  "2021-10-05 14:35:30.564429",
  "2021-10-05 14:35:30.564432",
  "2021-10-05 14:35:30.564434",
  "2021-10-05 14:35:30.564437",
  "2021-10-05 14:35:30.564439",
  "2021-10-05 14:35:30.564442",
  "2021-10-05 14:35:30.564444",
  "2021-10-05 14:35:30.564447",
  "2021-10-05 14:35:30.564450",
  "2021-10-05 14:35:30.564453",
  "2021-10-05 14:35:30.564455",
  "2021-10-05 14:35:30.564458",
  "2021-10-05 14:35:30.564460",
  "2021-10-05 14:35:30.564463",
  "2021-10-05 14:35:30.564465",
  "2021-10-05 14:35:30.564468",
  "2021-10-05 14:35:30.564470",
  "2021-10-05 14:35:30.564473",
  "2021-10-05 14:35:30.564475",
  "2021-10-05 14:35:30.564478",
  "2021-10-05 14:35:30.564480",
  "2021-10-05 14:35:30.564483",
  "2021-10-05 14:35:30.564485",
  "2021-10-05 14:35:30.564488",
  "2021-10-05 14:35:30.564490",
  "2021-10-05 14:35:30.564493",
  "2021-10-05 14:35:30.564495",
  "2021-10-05 14:35:30.564497",
  "2021-10-05 14:35:30.564500",
  "2021-10-05 14:35:30.564502",
  "2021-10-05 14:35:30.564505",
  "2021-10-05 14:35:30.564507",
  "2021-10-05 14:35:30.564510",
  "2021-10-05 14:35:30.564512",
  "2021-10-05 14:35:30.564515",
  "2021-10-05 14:35:30.564517",
  "2021-10-05 14:35:30.564520",
  "2021-10-05 14:35:30.564522",
  "2021-10-05 14:35:30.564525",
  "2021-10-05 14:35:30.564527",
  "2021-10-05 14:35:30.564530",
  "2021-10-05 14:35:30.564532",
  "2021-10-05 14:35:30.564534",
  "2021-10-05 14:35:30.564537",
  "2021-10-05 14:35:30.564540",
  "2021-10-05 14:35:30.564542",
  "2021-10-05 14:35:30.564545",
  "2021-10-05 14:35:30.564547",
  "2021-10-05 14:35:30.564550",
  "2021-10-05 14:35:30.564552",
  "2021-10-05 14:35:30.564555",
  "2021-10-05 14:35:30.564557",
  "2021-10-05 14:35:30.564560",
  "2021-10-05 14:35:30.564562",
  "2021-10-05 14:35:30.564565",
  "2021-10-05 14:35:30.564567",
  "2021-10-05 14:35:30.564570",
  "2021-10-05 14:35:30.564572",
  "2021-10-05 14:35:30.564575",
  "2021-10-05 14:35:30.564577",
  "2021-10-05 14:35:30.564579",
  "2021-10-05 14:35:30.564582",
  "2021-10-05 14:35:30.564584",
  "2021-10-05 14:35:30.564587",
  "2021-10-05 14:35:30.564589",
  "2021-10-05 14:35:30.564592",
  "2021-10-05 14:35:30.564594",
  "2021-10-05 14:35:30.564597",
  "2021-10-05 14:35:30.564599",
  "2021-10-05 14:35:30.564602",
  "2021-10-05 14:35:30.564604",
  "2021-10-05 14:35:30.564607",
  "2021-10-05 14:35:30.564609",
  "2021-10-05 14:35:30.564612",
  "2021-10-05 14:35:30.564614",
  "2021-10-05 14:35:30.564616",
  "2021-10-05 14:35:30.564619",
  "2021-10-05 14:35:30.564621",
  "2021-10-05 14:35:30.564624",
  "2021-10-05 14:35:30.564627",
  "2021-10-05 14:35:30.564629",
  "2021-10-05 14:35:30.564632",
  "2021-10-05 14:35:30.564634",
  "2021-10-05 14:35:30.564637",
  "2021-10-05 14:35:30.564641",
  "2021-10-05 14:35:30.564644",
  "2021-10-05 14:35:30.564646",
  "2021-10-05 14:35:30.564649",
  "2021-10-05 14:35:30.564651",
  "2021-10-05 14:35:30.564654",
  "2021-10-05 14:35:30.564656",
  "2021-10-05 14:35:30.564659",
  "2021-10-05 14:35:30.564661",
  "2021-10-05 14:35:30.564663",
  "2021-10-05 14:35:30.564666",
  "2021-10-05 14:35:30.564668",
  "2021-10-05 14:35:30.564671",
  "2021-10-05 14:35:30.564673",
  "2021-10-05 14:35:30.564676",
  "2021-10-05 14:35:30.564678",
  "2021-10-05 14:35:30.564680",
  "2021-10-05 14:35:30.564683",
  "2021-10-05 14:35:30.564685",
  "2021-10-05 14:35:30.564688",
  "2021-10-05 14:35:30.564690",
  "2021-10-05 14:35:30.564693",
  "2021-10-05 14:35:30.564695",
  "2021-10-05 14:35:30.564698",
  "2021-10-05 14:35:30.564700",
  "2021-10-05 14:35:30.564702",
  "2021-10-05 14:35:30.564705",
  "2021-10-05 14:35:30.564707",
  "2021-10-05 14:35:30.564710",
  "2021-10-05 14:35:30.564712",
  "2021-10-05 14:35:30.564715",
  "2021-10-05 14:35:30.564717",
  "2021-10-05 14:35:30.564720",
  "2021-10-05 14:35:30.564722",
  "2021-10-05 14:35:30.564725",
  "2021-10-05 14:35:30.564727",
  "2021-10-05 14:35:30.564729",
  "2021-10-05 14:35:30.564732",
  "2021-10-05 14:35:30.564734",
  "2021-10-05 14:35:30.564737",
  "2021-10-05 14:35:30.564739",
  "2021-10-05 14:35:30.564742",
  "2021-10-05 14:35:30.564744",
  "2021-10-05 14:35:30.564747",
  "2021-10-05 14:35:30.564751",
  "2021-10-05 14:35:30.564754",
  "2021-10-05 14:35:30.564756",
  "2021-10-05 14:35:30.564759",
  "2021-10-05 14:35:30.564761",
  "2021-10-05 14:35:30.564764",
  "2021-10-05 14:35:30.564766",
  "2021-10-05 14:35:30.564769",
  "2021-10-05 14:35:30.564771",
  "2021-10-05 14:35:30.564774",
  "2021-10-05 14:35:30.564776",
  "2021-10-05 14:35:30.564779",
  "2021-10-05 14:35:30.564781",
  "2021-10-05 14:35:30.564783",
  "2021-10-05 14:35:30.564786",
  "2021-10-05 14:35:30.564788",
  "2021-10-05 14:35:30.564791",
  "2021-10-05 14:35:30.564793",
  "2021-10-05 14:35:30.564796",
  "2021-10-05 14:35:30.564798",
  "2021-10-05 14:35:30.564801",
  "2021-10-05 14:35:30.564803",
  "2021-10-05 14:35:30.564806",
  "2021-10-05 14:35:30.564808",
  "2021-10-05 14:35:30.564811",
  "2021-10-05 14:35:30.564813",
  "2021-10-05 14:35:30.564816",
  "2021-10-05 14:35:30.564818",
  "2021-10-05 14:35:30.564821",
  "2021-10-05 14:35:30.564823",
  "2021-10-05 14:35:30.564826",
  "2021-10-05 14:35:30.564828",
  "2021-10-05 14:35:30.564831",
  "2021-10-05 14:35:30.564833",
  "2021-10-05 14:35:30.564836",
  "2021-10-05 14:35:30.564838",
  "2021-10-05 14:35:30.564841",
  "2021-10-05 14:35:30.564844",
  "2021-10-05 14:35:30.564846",
  "2021-10-05 14:35:30.564849",
  "2021-10-05 14:35:30.564851",
  "2021-10-05 14:35:30.564854",
  "2021-10-05 14:35:30.564856",
  "2021-10-05 14:35:30.564859",
  "2021-10-05 14:35:30.564862",
  "2021-10-05 14:35:30.564864",
  "2021-10-05 14:35:30.564867",
  "2021-10-05 14:35:30.564869"
)

for ($i = 0; $i -lt $newTimeTaken.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimeTaken[$i]
}

# --- Add the new "metadata" worksheet (after "data") ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

$headerRange = $metaSheet.Range("B1:G1")
$headerValues = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $cell = $metaSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headerValues[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$indexCell = $metaSheet.Cells.Item(2, 1)
$indexCell.Value = 0
$indexCell.Font.Bold = $true
$indexCell.Borders.LineStyle = 1
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160

$metaSheet.Range("B2").Value = "Regression"
$metaSheet.Range("C2").Value = 206
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.379"
$metaSheet.Range("E2").Value = "2021-10-04T04:28:36.769434Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:30.560650"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/206/?format=json"
